# Bulkupload with all fields: combine firstname/lastname into fullname and
# add the remaining CRM lead columns; refresh header row, column widths,
# sheet dimension and current selection to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order (A1:Y1) replacing the old 10-column (A1:J1) layout.
$headers = @(
    "leadowner",
    "source",
    "fullname",
    "contact",
    "company",
    "territory",
    "country",
    "requirements",
    "status",
    "primarycategory",
    "secondarycategory",
    "leadfor",
    "email",
    "whatsapp",
    "designation",
    "address",
    "region",
    "state",
    "isfca",
    "ivrticketcode",
    "isivrticketopen",
    "warrantystatus",
    "domesticorexport",
    "referredby",
    "referrefto"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# The first ten header cells (A1:J1) already carried the bold/italic header
# style; stamp the same style onto the newly added K1:Y1 header cells so the
# whole row 1 is formatted consistently.
$ws.Range("A1").Copy()
$ws.Range("K1:Y1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Auto-fit every header column so stored widths match the longer labels
# (mirrors Excel's bestFit recompute after editing the header row).
$ws.Range("A1:Y1").EntireColumn.AutoFit()

# Move the active selection to reflect where editing continued (column S,
# row 4) after the new columns were appended.
$ws.Range("S4").Select()
